$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "EU"
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = $false
$ws.Range("C5").Value = $false
$ws.Range("C6").Value = $true
$ws.Range("C7").Value = $false
$ws.Range("C8").Value = $false
$ws.Range("C9").Value = $true
$ws.Range("C10").Value = $false
$ws.Range("C11").Value = $true
$ws.Range("C12").Value = $true
$ws.Range("C13").Value = $true
$ws.Range("C14").Value = $true
$ws.Range("C15").Value = $false
$ws.Range("C16").Value = $true
$ws.Range("C17").Value = $true
$ws.Range("C18").Value = $true
$ws.Range("C19").Value = $false
$ws.Range("C20").Value = $true
$ws.Range("C21").Value = $true
$ws.Range("C22").Value = $false
$ws.Range("C23").Value = $true
$ws.Range("C24").Value = $true
$ws.Range("C25").Value = $false
$ws.Range("C26").Value = $false
$ws.Range("C27").Value = $false
$ws.Range("C28").Value = $true
$ws.Range("C29").Value = $false
$ws.Range("C30").Value = $true
$ws.Range("C31").Value = $false
$ws.Range("C32").Value = $true
$ws.Range("C33").Value = $false
$ws.Range("C34").Value = $false
$ws.Range("C35").Value = $true
$ws.Range("C36").Value = $true
$ws.Range("C37").Value = $false
$ws.Range("C38").Value = $false
$ws.Range("C39").Value = $false
$ws.Range("C40").Value = $true
$ws.Range("C41").Value = $true
$ws.Range("C42").Value = $true
$ws.Range("C43").Value = $false
$ws.Range("C44").Value = $false
$ws.Range("C45").Value = $false
$ws.Range("C46").Value = $true
$ws.Range("C47").Value = $false
$ws.Range("C48").Value = $true
$ws.Range("C49").Value = $true
$ws.Range("C50").Value = $true
$ws.Range("C51").Value = $false
$ws.Range("C52").Value = $true
$ws.Range("C50").Select()
